$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This price-history sheet stores every value (dates, prices, discounts,
# flags) as plain text (shared strings), not as native numbers/dates.
# Assigning .Value directly would let Excel "smart type" the new date
# ("2026-02-07") into a date serial and the numeric-looking strings into
# real numbers, plus it would allocate new number-format styles.
#
# Work around that by writing each value as a literal-string formula
# (=\"...\") first - this evaluates to a plain text result without any
# reinterpretation - then flattening the row to static values via
# copy / paste-special-values. That converts the formulas to literal
# text cells (shared strings) while leaving cell styling untouched,
# exactly mirroring the rest of the sheet.

$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1  # xlUp

$newDate = "2026-02-07"
$price = "113000"
$discount = "0"
$incredible = "0"

$ws.Cells.Item($newRow, 1).Formula = '="' + $newDate + '"'
$ws.Cells.Item($newRow, 2).Formula = '="' + $price + '"'
$ws.Cells.Item($newRow, 3).Formula = '="' + $discount + '"'
$ws.Cells.Item($newRow, 4).Formula = '="' + $incredible + '"'

$rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4))
$rowRange.Copy()
$rowRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0
